$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): extend with two more header values, matching style of O1 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the formatting (bold/border/centered) from the existing header cell O1
# onto the two new header cells so they match the rest of the header row.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Re-assert values (PasteSpecial only touches formatting, but ensure they are correct)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: flip values in columns I, K, M, O and add new columns P, Q ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P -> 2 (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q -> 2 (new)
}
